$d = $word.ActiveDocument

# 1) "Summoning Cthulhu" paragraph: drop the "inspired by Mara Averick / this tweet" lead-in,
#    keep only the sentence about the marine creature modification.
$d.Content.Find.Execute(
    "The name is inspired in an answer from Mara Averick to this tweet. It is a modification of the marine creature in polar coordinates:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "It is a modification of the marine creature in polar coordinates:", 2) | Out-Null

# 2) "Silk Knitting" paragraph: drop the "It is inspired by this other pattern" lead-in,
#    keep only the ". A lot of " -> "A lot of " trim; leave the following runs untouched.
$d.Content.Find.Execute(
    "It is inspired by this other pattern. A lot of ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A lot of ", 2) | Out-Null
